$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Cells whose new Price text would parse as a plain number need the cell
# forced to Text format first so Excel keeps storing the literal string
# (matching the source site's formatting, e.g. trailing zeros) instead of
# silently converting it to a numeric value.

$ws.Range("D2").Value = "47.270.00"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "2.492.52"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.75"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.42"
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.21"
$ws.Range("E10").Value = "  +5.08%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.14"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").Value = "2.881.22"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "2.490.88"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.845"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "47.177.00"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.17"
$ws.Range("E19").Value = "  +3.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.63"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").Value = "0.0₃0940"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("E22").Value = "  +14.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.42"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "245.55"
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.72"
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("E30").Value = "  +3.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.64"
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.78"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.14"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.34"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0784"
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.75"
$ws.Range("E37").Value = "  +2.46%  "
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.01"
$ws.Range("E40").Value = "  +7.62%  "
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.22"
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "1.997.98"
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("E46").Value = "  +2.04%  "
$ws.Range("E47").Value = "  -4.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.15"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.08"
$ws.Range("E50").Value = "  -5.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.44"
$ws.Range("E51").Value = "  +3.46%  "
